$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase the header row (row 1) labels, B1:X1 (A1 "NIK" and C1 "KTP" are
# already upper-case so they are left untouched, matching the source diff).
$ws.Range("B1").Value = "NAME"
$ws.Range("D1").Value = "PHONE"
$ws.Range("E1").Value = "EMAIL"
$ws.Range("F1").Value = "TIMEZONE"
$ws.Range("G1").Value = "REKENING"
$ws.Range("H1").Value = "BANK"
$ws.Range("I1").Value = "JOIN DATE"
$ws.Range("J1").Value = "AGENCY"
$ws.Range("K1").Value = "GENDER"
$ws.Range("L1").Value = "BIRTH DATE"
$ws.Range("M1").Value = "POSITION"
$ws.Range("N1").Value = "STATUS"
$ws.Range("O1").Value = "STORE"
$ws.Range("P1").Value = "SUBAREA"
$ws.Range("Q1").Value = "AREA"
$ws.Range("R1").Value = "REGION"
$ws.Range("S1").Value = "ACCOUNT"
$ws.Range("T1").Value = "CHANNEL"
$ws.Range("U1").Value = "PASSWORD"
$ws.Range("V1").Value = "SALES TIER"
$ws.Range("W1").Value = "TIMEZONE STORE"
$ws.Range("X1").Value = "EDUCATION"

# Bold the (empty) trailing header cell Y1.
$ws.Range("Y1").Font.Bold = $true

# Scroll the viewport over and move the selection to Y1, matching the
# author's final cursor position after editing the last header cell.
$excel.ActiveWindow.ScrollColumn = 18
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("Y1").Select()
